# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# The match-data rows had their per-match "id" (column B) and all of the
# associated odds/result columns (E through AD) mixed up between several
# rows. This script fixes the mapping by moving each row's B:AD content to
# the row it actually belongs to:
#   - rows 27 and 28 are swapped with each other
#   - rows 31, 32, 33, 36 are rotated (31<-32, 32<-33, 33<-36, 36<-31)
#   - rows 197 and 201 are swapped with each other
# Column A (the running counter) stays put on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

function Get-RowValues($row) {
    $values = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $values += ,($ws.Cells.Item($row, $c).Value())
    }
    return $values
}

function Set-RowValues($row, $values) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value() = $values[$i]
        $i = $i + 1
    }
}

# Capture the current ("before") content of every affected row first, so
# that writes to one row never clobber data we still need to read from
# another row.
$row27 = Get-RowValues 27
$row28 = Get-RowValues 28

$row31 = Get-RowValues 31
$row32 = Get-RowValues 32
$row33 = Get-RowValues 33
$row36 = Get-RowValues 36

$row197 = Get-RowValues 197
$row201 = Get-RowValues 201

# Swap 27 <-> 28
Set-RowValues 27 $row28
Set-RowValues 28 $row27

# Rotate 31 <- 32 <- 33 <- 36 <- 31
Set-RowValues 31 $row32
Set-RowValues 32 $row33
Set-RowValues 33 $row36
Set-RowValues 36 $row31

# Swap 197 <-> 201
Set-RowValues 197 $row201
Set-RowValues 201 $row197
